$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rodada 14")

# Clear the Mandante_Pontos (D) and Visitante_Pontos (F) values for rows 2-5
$ws.Range("D2:D5").ClearContents()
$ws.Range("F2:F5").ClearContents()
